$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5: Date (A), Start (B), Stop (C) - reuse existing number formats
# from row 2 via copy/paste-special so no new custom style gets created.
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = 42375

$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = 0.58680555555555558

$ws.Range("C2").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = 0.59722222222222221

# Notes column for the new row
$ws.Range("E5").Value = "Added pay period formatting."

# Move the active selection to the newly added last cell, matching
# the updated sheetView selection in the workbook.
$ws.Range("E5").Select()
